$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-unused rows 8-20 first
$ws.Rows("8:20").Delete()

# Update consolidated token rows
$ws.Range("A2").Value = "('Goblin', ['Token Creature — Goblin', '1/1'])"
$ws.Range("A3").Value = "('Pegasus', ['Token Creature — Pegasus', 'Flying', '1/1'])"
$ws.Range("A4").Value = "('Sheep', ['Token Creature — Sheep', '2/2'])"
$ws.Range("A5").Value = "('Soldier', ['Token Creature — Soldier', '1/1'])"
$ws.Range("A6").Value = "('Squirrel', ['Token Creature — Squirrel', '1/1'])"
$ws.Range("A7").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"
